# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.298.95"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.928.04"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.92"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'0.7182"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.3202"
$ws.Range("E8").Value = "  -4.79%  "
$ws.Range("D9").Value = "'27.85"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("D10").Value = "'0.07083"
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").Value = "'0.07985"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'1.930.85"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'5.376"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "'94.72"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'14.64"
$ws.Range("D17").Value = "'30.294.10"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'257.20"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'0.000008092"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "'5.756"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'2.184.23"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'6.821"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("D25").Value = "'9.536"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "'164.48"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "'19.09"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").Value = "'2.271"
$ws.Range("E28").Value = "  -7.54%  "
$ws.Range("D29").Value = "'0.1269"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("D31").Value = "'1.530"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "'4.403"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").Value = "'4.137"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").Value = "'0.05137"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "'0.7448"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'0.01984"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "'2.799"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "'78.15"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("D41").Value = "'6.374"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("D42").Value = "'0.4507"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "'1.996"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "'0.8462"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'0.9997"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'100.70"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "'9.795"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'36.85"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "'953.88"
$ws.Range("E50").Value = "  +8.13%  "
$ws.Range("D51").Value = "'0.4216"
$ws.Range("E51").Value = "  +0.10%  "
